$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column J width (was 15.8 -> ~36.52). Excel's ColumnWidth setter snaps to
# whole-pixel increments, so we pick the value that lands closest to the
# target serialized width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 35.67

# ---------------------------------------------------------------------------
# Row height tweaks (content re-wrap after the data fixes below).
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 13.8
$ws.Rows.Item(11).RowHeight = 13.8
$ws.Rows.Item(13).RowHeight = 13.8
$ws.Rows.Item(15).RowHeight = 23.95

# ---------------------------------------------------------------------------
# Row 15 ("TC R_14"): debugged with a fresh serial key / pos user / the
# multi-item buy basket (rich text: three SKUs, middle one in a second
# font, the separator left uncoloured).
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "307260624Wa9 "
$ws.Range("B15").Font.Name = "Times New Roman"
$ws.Range("B15").Font.Size = 10
$ws.Range("B15").Font.Color = 0

$ws.Range("E15").Value = "userone_p10"

$j15 = $ws.Range("J15")
$j15.Value = "8906118410781 : 1, 8906118412556 : 1, 8906118412662:1"
$j15.Font.Name = "Arial"
$j15.Font.Size = 10
$j15.Font.Color = 0
$j15.Characters(20, 13).Font.Name = "Arial"
$j15.Characters(20, 13).Font.Size = 10
$j15.Characters(20, 13).Font.Color = 0
$j15.Characters(33, 6).Font.Name = "Arial"
$j15.Characters(33, 6).Font.Size = 10
$j15.Characters(39, 15).Font.Name = "Arial"
$j15.Characters(39, 15).Font.Size = 10
$j15.Characters(39, 15).Font.Color = 0

# ---------------------------------------------------------------------------
# Row 16 ("TC R_15"): reverts back to the common serial key / pos user
# (it used to hold the now-relocated custom data).
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "307260624P3E"
$ws.Range("B16").Font.Name = "Arial"
$ws.Range("B16").Font.Size = 10
$ws.Range("B16").Font.Color = 0

$ws.Range("E16").Value = "userone_p1 "

# ---------------------------------------------------------------------------
# New rows 21-35 (TC R_20 .. R_34), identical template to rows 17-20 minus
# the trailing Q "Return" marker.
# ---------------------------------------------------------------------------
for ($i = 21; $i -le 35; $i++) {
  $ws.Rows.Item($i).RowHeight = 13.8

  $ws.Cells.Item($i, 1).Value = "R_" + ($i - 1)
  $ws.Cells.Item($i, 1).Style = $ws.Cells.Item(20, 1).Style

  $ws.Cells.Item($i, 2).Value = "userone_p1 "
  $ws.Cells.Item($i, 2).Style = $ws.Cells.Item(20, 2).Style

  $ws.Cells.Item($i, 3).Value = "Index9QA"
  $ws.Cells.Item($i, 3).Style = $ws.Cells.Item(20, 3).Style

  $ws.Cells.Item($i, 4).Value = 123456
  $ws.Cells.Item($i, 4).Style = $ws.Cells.Item(20, 4).Style

  $ws.Cells.Item($i, 5).Value = "userone_p1 "
  $ws.Cells.Item($i, 5).Style = $ws.Cells.Item(20, 5).Style

  $ws.Cells.Item($i, 6).Value = 123456
  $ws.Cells.Item($i, 6).Style = $ws.Cells.Item(20, 6).Style

  $ws.Cells.Item($i, 7).Value = 1000
  $ws.Cells.Item($i, 7).Style = $ws.Cells.Item(20, 7).Style

  $ws.Cells.Item($i, 8).Value = 400
  $ws.Cells.Item($i, 8).Style = $ws.Cells.Item(20, 8).Style

  $ws.Cells.Item($i, 9).Value = "NULL"
  $ws.Cells.Item($i, 9).Style = $ws.Cells.Item(20, 9).Style

  $ws.Cells.Item($i, 10).Value = "8906118410781 : 1"
  $ws.Cells.Item($i, 10).Style = $ws.Cells.Item(20, 10).Style

  $ws.Cells.Item($i, 11).Value = "NULL"
  $ws.Cells.Item($i, 11).Style = $ws.Cells.Item(20, 11).Style

  $ws.Cells.Item($i, 12).Value = 45384
  $ws.Cells.Item($i, 12).Style = $ws.Cells.Item(20, 12).Style

  $ws.Cells.Item($i, 13).Value = "Maharashtra"
  $ws.Cells.Item($i, 13).Style = $ws.Cells.Item(20, 13).Style

  $ws.Cells.Item($i, 14).Value = "Pune"
  $ws.Cells.Item($i, 14).Style = $ws.Cells.Item(20, 14).Style

  $ws.Cells.Item($i, 15).Value = "Dummy"
  $ws.Cells.Item($i, 15).Style = $ws.Cells.Item(20, 15).Style

  $ws.Cells.Item($i, 16).Value = "NULL"
  $ws.Cells.Item($i, 16).Style = $ws.Cells.Item(20, 16).Style
}

# ---------------------------------------------------------------------------
# Selection left where the author was last looking.
# ---------------------------------------------------------------------------
$ws.Range("B19").Select()
